$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry per-record data for rows 73-77.
$cols = @("A","B","D","E","F","G","H","Q","R")

# Snapshot current values for rows 73-77 before overwriting anything,
# since the update cyclically shifts the records between rows.
$rows = @(73,74,75,76,77)
$data = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $data[$r] = $rowVals
}

# Record at row 77 moves up to row 73; records at rows 73-76 each shift
# down by one row (74->75, 75->76, 76->77, 73->74).
$mapping = @{
    73 = 77
    74 = 73
    75 = 74
    76 = 75
    77 = 76
}

foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    $srcVals = $data[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
